# Updating numbers for calcs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EEU data")

# --- Row 2 ---
$ws.Range("H2").Value = 3.2380660967337498
$ws.Range("J2").Value = 0.12870000000000001

# --- Row 3 ---
$ws.Range("H3").Value = 3.2380660967337498
$ws.Range("J3").Value = 2.63
$ws.Range("Q3").Value = 27897.989870000001
$ws.Range("R3").Value = 33446.43
$ws.Range("T3").Value = 33037.919999999998
$ws.Range("V3").NumberFormat = "#,##0.00000"
$ws.Range("V3").Value = 2730.84854
$ws.Range("W3").NumberFormat = "#,##0.00000"
$ws.Range("W3").Value = 2709.6940300000001

# --- Row 4 ---
$ws.Range("H4").Value = 3.2380660967337498
$ws.Range("J4").Value = 0.12870000000000001
$ws.Range("Q4").Value = 27897.989870000001

# --- Row 5 ---
$ws.Range("H5").Value = 3.2380660967337498
$ws.Range("Q5").Value = 27897.989870000001
$ws.Range("R5").Value = 33446.43
$ws.Range("T5").Value = 33037.919999999998
$ws.Range("V5").NumberFormat = "#,##0.00000"
$ws.Range("V5").Value = 2730.84854
$ws.Range("W5").NumberFormat = "#,##0.00000"
$ws.Range("W5").Value = 2709.6940300000001

# --- Row 6 ---
$ws.Range("H6").Value = 3.2380660967337498
$ws.Range("Q6").Value = 27897.989870000001

# --- Selection moved to AB3 ---
$null = $ws.Range("AB3").Select()
